{"js": "// Update the date line and the 25 division-fact answers in the practice\n// table. Both the date paragraph and every populated table cell contain\n// exactly one paragraph of text, so we can walk `body.paragraphs` in\n// document order (this also covers the text living inside table cells)\n// and overwrite each non-empty paragraph's text with its replacement,\n// preserving the existing run formatting via `insertText(..., \"Replace\")`.\n\nconst replacements = [\n  \"2025-12-03 Wednesday\",\n  \"85\u00f76=14, 1\",\n  \"48\u00f74=12, 0\",\n  \"84\u00f72=42, 0\",\n  \"80\u00f73=26, 2\",\n  \"10\u00f76=1, 4\",\n  \"31\u00f76=5, 1\",\n  \"92\u00f73=30, 2\",\n  \"25\u00f74=6, 1\",\n  \"74\u00f74=18, 2\",\n  \"50\u00f75=10, 0\",\n  \"69\u00f74=17, 1\",\n  \"59\u00f78=7, 3\",\n  \"45\u00f74=11, 1\",\n  \"21\u00f77=3, 0\",\n  \"26\u00f73=8, 2\",\n  \"46\u00f76=7, 4\",\n  \"99\u00f79=11, 0\",\n  \"22\u00f72=11, 0\",\n  \"35\u00f77=5, 0\",\n  \"71\u00f72=35, 1\",\n  \"34\u00f73=11, 1\",\n  \"66\u00f72=33, 0\",\n  \"35\u00f72=17, 1\",\n  \"21\u00f73=7, 0\",\n  \"81\u00f74=20, 1\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idx = 0;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === \"\") {\n    continue;\n  }\n  if (idx >= replacements.length) {\n    break;\n  }\n  paragraph.insertText(replacements[idx], Word.InsertLocation.replace);\n  idx++;\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-fact answers in the practice\n# table. The table has a duplicate answer string (\"59\u00f76=9, 5\") appearing\n# in two different cells that map to two different replacements, so a\n# single document-wide Find/Replace-All would be ambiguous. Instead we\n# scope each Find to the specific table cell's Range and use\n# wdReplaceOne (1) so only that cell's occurrence is touched.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph.\n$d.Paragraphs.Item(1).Range.Text = \"2025-12-03 Wednesday\"\n\n# 2) Division-fact table (single table in the document).\n$tbl = $d.Tables.Item(1)\n\nfunction Set-CellText($table, $row, $col, $oldText, $newText) {\n    $cell = $table.Cell($row, $col)\n    $range = $cell.Range\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    # wdFindContinue=1, wdReplaceOne=1 -> replace just this single match,\n    # scoped to the cell's own range.\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1) | Out-Null\n}\n\n# Row 1 (table row index 1)\nSet-CellText $tbl 1 1 \"53\u00f73=17, 2\" \"85\u00f76=14, 1\"\nSet-CellText $tbl 1 2 \"96\u00f74=24, 0\" \"48\u00f74=12, 0\"\nSet-CellText $tbl 1 3 \"30\u00f78=3, 6\" \"84\u00f72=42, 0\"\nSet-CellText $tbl 1 4 \"27\u00f72=13, 1\" \"80\u00f73=26, 2\"\nSet-CellText $tbl 1 5 \"59\u00f76=9, 5\" \"10\u00f76=1, 4\"\n\n# Row 5 (table row index 5)\nSet-CellText $tbl 5 1 \"57\u00f78=7, 1\" \"31\u00f76=5, 1\"\nSet-CellText $tbl 5 2 \"59\u00f76=9, 5\" \"92\u00f73=30, 2\"\nSet-CellText $tbl 5 3 \"49\u00f78=6, 1\" \"25\u00f74=6, 1\"\nSet-CellText $tbl 5 4 \"37\u00f73=12, 1\" \"74\u00f74=18, 2\"\nSet-CellText $tbl 5 5 \"81\u00f75=16, 1\" \"50\u00f75=10, 0\"\n\n# Row 9 (table row index 9)\nSet-CellText $tbl 9 1 \"61\u00f74=15, 1\" \"69\u00f74=17, 1\"\nSet-CellText $tbl 9 2 \"50\u00f78=6, 2\" \"59\u00f78=7, 3\"\nSet-CellText $tbl 9 3 \"42\u00f72=21, 0\" \"45\u00f74=11, 1\"\nSet-CellText $tbl 9 4 \"73\u00f79=8, 1\" \"21\u00f77=3, 0\"\nSet-CellText $tbl 9 5 \"58\u00f74=14, 2\" \"26\u00f73=8, 2\"\n\n# Row 13 (table row index 13)\nSet-CellText $tbl 13 1 \"96\u00f77=13, 5\" \"46\u00f76=7, 4\"\nSet-CellText $tbl 13 2 \"44\u00f73=14, 2\" \"99\u00f79=11, 0\"\nSet-CellText $tbl 13 3 \"22\u00f73=7, 1\" \"22\u00f72=11, 0\"\nSet-CellText $tbl 13 4 \"47\u00f75=9, 2\" \"35\u00f77=5, 0\"\nSet-CellText $tbl 13 5 \"29\u00f79=3, 2\" \"71\u00f72=35, 1\"\n\n# Row 17 (table row index 17)\nSet-CellText $tbl 17 1 \"89\u00f79=9, 8\" \"34\u00f73=11, 1\"\nSet-CellText $tbl 17 2 \"98\u00f78=12, 2\" \"66\u00f72=33, 0\"\nSet-CellText $tbl 17 3 \"57\u00f75=11, 2\" \"35\u00f72=17, 1\"\nSet-CellText $tbl 17 4 \"15\u00f77=2, 1\" \"21\u00f73=7, 0\"\nSet-CellText $tbl 17 5 \"63\u00f79=7, 0\" \"81\u00f74=20, 1\"\n"}
